# "Create Deal Code Update"
# Adds the CreateDeal_/ShareDeal_ automation test rows to the Login sheet,
# re-selects the Login tab (it was previously on "Shipper Admin"), and
# updates the active cell/selection to match the newly appended data.

$wb = $excel.ActiveWorkbook
$login = $wb.Worksheets.Item("Login")

# New automation test rows (Automation Test ID, UserName, Password, Status)
$rows = @(
    @("CreateDeal_TC001", "rogerdeals21+stan@gmail.com", "arewethere?", "Login successful"),
    @("CreateDeal_TC002", "rogerdeals21+rick@gmail.com", "arewethere?", "Login successful"),
    @("CreateDeal_TC003", "rogerdeals21+john@gmail.com", "arewethere?", "Login successful"),
    @("ShareDeal_TC001",  "rogerdeals21+stan@gmail.com", "arewethere?", "Login successful"),
    @("ShareDeal_TC001(2)", "rogerdeals21+john@gmail.com", "arewethere?", "Login successful"),
    @("ShareDeal_TC002",  "rogerdeals21+rick@gmail.com", "arewethere?", "Login successful"),
    @("ShareDeal_TC002(2)", "rogerdeals21+john@gmail.com", "arewethere?", "Login successful")
)

$startRow = 19
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $cellA = $login.Cells.Item($r, 1)
    $cellA.Value = $data[0]
    $cellA.VerticalAlignment = -4108   # xlCenter - matches the new style used for column A

    $login.Cells.Item($r, 2).Value = $data[1]
    $login.Cells.Item($r, 3).Value = $data[2]
    $login.Cells.Item($r, 4).Value = $data[3]
}

# The Login sheet becomes the active tab/sheet again (it had moved to
# "Shipper Admin" before this change).
$login.Activate()

# Update the visible selection to the last newly-added row.
$login.Range("B24").Select()

# "Shipper Admin" is no longer the selected tab.
$shipperAdmin = $wb.Worksheets.Item("Shipper Admin")
